# Update the cryptocurrency price list (Price and Volume(1h) columns)
# with freshly scraped values, preserving the existing text-cell storage
# (prices like "524.50" or "2.606.66" must remain literal text, not be
# auto-converted to numbers by Excel's type inference).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to keep/receive a literal text value, regardless of
    # whether the string looks like a number, then drop the temporary
    # "Text" number format so no extra style is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$updates = @(
    @{ Row = 2;  D = "61.022.61";  E = "  +0.84%  " }
    @{ Row = 3;  D = "2.606.66";   E = "  +0.65%  " }
    @{ Row = 4;  E = "  +0.15%  " }
    @{ Row = 5;  D = "524.50";     E = "  +3.23%  " }
    @{ Row = 6;  D = "155.17";     E = "  +1.05%  " }
    @{ Row = 7;  D = "0.998";      E = "  +0.08%  " }
    @{ Row = 8;  E = "  +2.02%  " }
    @{ Row = 9;  E = "  +2.41%  " }
    @{ Row = 10; D = "0.106";      E = "  +1.64%  " }
    @{ Row = 11; D = "0.348";      E = "  -0.21%  " }
    @{ Row = 12; E = "  +1.20%  " }
    @{ Row = 13; D = "3.060.20";   E = "  +0.96%  " }
    @{ Row = 14; D = "60.979.60";  E = "  +0.92%  " }
    @{ Row = 15; D = "21.74";      E = "  +0.55%  " }
    @{ Row = 16; D = "0.0000142";  E = "  +0.95%  " }
    @{ Row = 17; D = "2.604.72";   E = "  +1.03%  " }
    @{ Row = 18; D = "4.77";       E = "  -0.23%  " }
    @{ Row = 19; D = "355.78";     E = "  +2.74%  " }
    @{ Row = 20; D = "10.59";      E = "  +1.45%  " }
    @{ Row = 21; E = "  +2.11%  " }
    @{ Row = 22; E = "  -0.10%  " }
    @{ Row = 23; D = "61.02";      E = "  +2.06%  " }
    @{ Row = 24; E = "  +1.65%  " }
    @{ Row = 25; E = "  +0.60%  " }
    @{ Row = 26; D = "2.719.11";   E = "  +1.30%  " }
    @{ Row = 27; D = "0.998";      E = "  +0.15%  " }
    @{ Row = 28; D = "0.0₃0850";   E = "  +0.67%  " }
    @{ Row = 29; E = "  +0.39%  " }
    @{ Row = 30; E = "  +0.06%  " }
    @{ Row = 31; D = "6.30";       E = "  +10.15%  " }
    @{ Row = 32; D = "19.42";      E = "  +0.38%  " }
    @{ Row = 33; D = "1.61";       E = "  +2.92%  " }
    @{ Row = 34; D = "148.20";     E = "  -3.41%  " }
    @{ Row = 35; D = "4.21";       E = "  +5.33%  " }
    @{ Row = 36; D = "1.21";       E = "  +1.42%  " }
    @{ Row = 37; D = "0.919";      E = "  +8.22%  " }
    @{ Row = 38; D = "0.889";      E = "  +4.39%  " }
    @{ Row = 39; E = "  +1.87%  " }
    @{ Row = 40; D = "3.81";       E = "  +1.74%  " }
    @{ Row = 41; E = "  +1.15%  " }
    @{ Row = 42; D = "290.95";     E = "  -2.39%  " }
    @{ Row = 43; E = "  +2.21%  " }
    @{ Row = 44; E = "  +0.33%  " }
    @{ Row = 45; D = "0.0561";     E = "  +0.74%  " }
    @{ Row = 46; D = "0.997";      E = "  +0.18%  " }
    @{ Row = 47; D = "19.57";      E = "  -1.19%  " }
    @{ Row = 48; D = "4.93";       E = "  +0.42%  " }
    @{ Row = 49; E = "  +2.27%  " }
    @{ Row = 50; E = "  +0.22%  " }
    @{ Row = 51; D = "19.26";      E = "  +8.92%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("D")) {
        # Price column values such as "524.50" or "2.606.66" look numeric
        # (or partially numeric) to Excel's auto-detection, so they need
        # to be forced to text explicitly.
        Set-TextValue $ws.Range("D$row") $u.D
    }
    # Volume(1h) values always contain '%' and padding spaces, so they are
    # never misread as numbers and can be assigned directly.
    $ws.Range("E$row").Value = $u.E
}
